# daily auto push: 2026-01-05 02:20 UTC
# Append 26 new data rows (2026/01/01 - 2026/01/05) below the existing
# table, which currently runs through row 536 (2025/12/31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(537, "2026/01/01", "木", 2, 129),
  @(538, "2026/01/01", "木", 5, 119),
  @(539, "2026/01/01", "木", 13, 133),
  @(540, "2026/01/01", "木", 16, 109),
  @(541, "2026/01/01", "木", 19, 120),
  @(542, "2026/01/02", "金", 1, 105),
  @(543, "2026/01/02", "金", 5, 109),
  @(544, "2026/01/02", "金", 8, 110),
  @(545, "2026/01/02", "金", 13, 132),
  @(546, "2026/01/02", "金", 16, 145),
  @(547, "2026/01/02", "金", 19, 157),
  @(548, "2026/01/02", "金", 22, 165),
  @(549, "2026/01/03", "土", 1, 174),
  @(550, "2026/01/03", "土", 4, 192),
  @(551, "2026/01/03", "土", 7, 189),
  @(552, "2026/01/03", "土", 13, 201),
  @(553, "2026/01/03", "土", 16, 201),
  @(554, "2026/01/03", "土", 19, 201),
  @(555, "2026/01/03", "土", 22, 194),
  @(556, "2026/01/04", "日", 2, 164),
  @(557, "2026/01/04", "日", 5, 166),
  @(558, "2026/01/04", "日", 7, 168),
  @(559, "2026/01/04", "日", 13, 173),
  @(560, "2026/01/04", "日", 22, 127),
  @(561, "2026/01/05", "月", 1, 118),
  @(562, "2026/01/05", "月", 7, 127)
)

foreach ($row in $rows) {
  $r = $row[0]

  # Column A holds a date-like string ("2026/01/01") that must stay plain
  # text (matches the rest of the sheet, which stores dates as text, not
  # as real Excel date serials). Mark the cell as Text first so Excel's
  # COM layer doesn't auto-convert the string into a date value, then
  # restore the default (unstyled) cell style so the new rows end up
  # without any explicit style, matching the rest of the sheet.
  $cellA = $ws.Cells.Item($r, 1)
  $cellA.NumberFormat = "@"
  $cellA.Value = $row[1]
  $cellA.Style = "Normal"

  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
}
